$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header style tweaks (A3/E3 and H5/H24 kept same semantic style; J10 keeps its numeric format) ---
# (styles already correct pre-existing; only new content below needs explicit styling)

# --- Row 27: nbmodels header now just "nbmodels" / 2, drop old rho/vp/vs mini headers in I27:K27 ---
$ws.Range("I27:K27").ClearContents()
$ws.Range("C27").Value = 2

# --- Clear the old M1-air / M2-clay / M3_water material table (rows 28-30, cols H:K) ---
$ws.Range("H28:K30").ClearContents()

# --- New model-table header row (row 30), columns H through P ---
$ws.Range("H30").Value = "domain_id"
$ws.Range("I30").Value = "material_id"
$ws.Range("J30").Value = "dens"
$ws.Range("K30").Value = "vp"
$ws.Range("L30").Value = "vs"
$ws.Range("M30").Value = "qkappa"
$ws.Range("N30").Value = "qmu"
$ws.Range("O30").Value = "ani"
$ws.Range("P30").Value = "# comment"
$ws.Range("H30:P30").Font.Bold = $true

# --- Row 31: material 1 / domain 2, clay ---
$ws.Range("H31").Value = 2
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = "1784.d0"
$ws.Range("K31").Value = "1700.d0"
$ws.Range("L31").Value = "500.0d0"
$ws.Range("M31").Value = 9999
$ws.Range("N31").Value = 9999
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = "# clay with specific weight 17.5 kN/m3"
$ws.Range("A31").Formula = '=TEXTJOIN(" ",TRUE,H31:P31)'

# --- Row 32: material 2 / domain 2, sand ---
$ws.Range("H32").Value = 2
$ws.Range("I32").Value = 2
$ws.Range("J32").Value = "1886.d0"
$ws.Range("K32").Value = "1600.d0"
$ws.Range("L32").Value = "400.0d0"
$ws.Range("M32").Value = 9999
$ws.Range("N32").Value = 9999
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = "# sand with specific weight 18.5 kN/m3"
$ws.Range("A32").Formula = '=TEXTJOIN(" ",TRUE,H32:P32)'

# --- Row 33: material 3 / domain 1, water ---
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = 3
$ws.Range("J33").Value = "1000.d0"
$ws.Range("K33").Value = "1480.d0"
$ws.Range("L33").Value = "0.000d0"
$ws.Range("M33").Value = 9999
$ws.Range("N33").Value = 9999
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = "# water burrow entrance"
$ws.Range("A33").Formula = '=TEXTJOIN(" ",TRUE,H33:P33)'

# --- Row 34: material 4 / domain 1, air ---
$ws.Range("H34").Value = 1
$ws.Range("I34").Value = 4
$ws.Range("J34").Value = "1.000d0"
$ws.Range("K34").Value = "343.0d0"
$ws.Range("L34").Value = "0.000d0"
$ws.Range("M34").Value = 9999
$ws.Range("N34").Value = 9999
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = "# air in burrow"
$ws.Range("A34").Formula = '=TEXTJOIN(" ",TRUE,H34:P34)'

# --- New empty, bold-styled cell next to the PML/BCs block ---
$ws.Range("G26").Font.Bold = $true

# --- I19 loses its (redundant) explicit number-format style ---
$ws.Range("I19").Style = "Normal"

# --- Column widths for the new columns (target XML widths 34.33203125 / 13.33203125) ---
$ws.Columns.Item(7).ColumnWidth = 33.498697916666664
$ws.Columns.Item(8).ColumnWidth = 12.498697916666666

# --- Selection moves to G23 ---
$ws.Range("G23").Select()
